# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2210"
#   "<header>_new" -> "<header>_FV2304"
# then format the data range as an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells in row 1 (A1:K1 = "_old" -> "_FV2210", L1:U1 = "_new" -> "_FV2304") ---
$oldToFV2210 = @{
    "A1" = "Segmentname_FV2210"
    "B1" = "Segmentgruppe_FV2210"
    "C1" = "Segment_FV2210"
    "D1" = "Datenelement_FV2210"
    "E1" = "Segment ID_FV2210"
    "F1" = "Code_FV2210"
    "G1" = "Qualifier_FV2210"
    "H1" = "Beschreibung_FV2210"
    "I1" = "Bedingungsausdruck_FV2210"
    "J1" = "Bedingung_FV2210"
}

$newToFV2304 = @{
    "L1" = "Segmentname_FV2304"
    "M1" = "Segmentgruppe_FV2304"
    "N1" = "Segment_FV2304"
    "O1" = "Datenelement_FV2304"
    "P1" = "Segment ID_FV2304"
    "Q1" = "Code_FV2304"
    "R1" = "Qualifier_FV2304"
    "S1" = "Beschreibung_FV2304"
    "T1" = "Bedingungsausdruck_FV2304"
    "U1" = "Bedingung_FV2304"
}

foreach ($addr in $oldToFV2210.Keys) {
    $ws.Range($addr).Value = $oldToFV2210[$addr]
}
foreach ($addr in $newToFV2304.Keys) {
    $ws.Range($addr).Value = $newToFV2304[$addr]
}

# --- 2. Turn the data range into a proper Excel Table (ListObject) ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U86"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze the header row (top row) ---
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
